# Edit script for "Cannabis: Cash Crop or Crime Boom?" presentation
# Applies 4 text changes described by the commit diff:
#   1. Slide 3 body, paragraph 3: append " with \u03b1 = 0.05" before trailing colon
#   2. Slide 4 body, paragraph 2: "not significant" -> "not statistically significant"
#   3. Slide 8 title: split "Results and Next Steps" into two runs: "Results " + "and Next Steps"
#   4. Slide 8 body, paragraph 1: reworded conclusion sentence
#
# NOTE: when a TextRange.Text assignment shares a long common prefix/suffix with
# the text it replaces, the host splits the paragraph into multiple runs around
# the differing span (it keeps the unchanged prefix/suffix as separate runs).
# Since the target XML keeps each of these paragraphs as a single run, we first
# set the paragraph text to an unrelated placeholder (breaking prefix/suffix
# overlap) and then set it to the final desired text so the whole paragraph
# collapses back into one run with the same default rPr.

$p = $ppt.ActivePresentation

# --- 1. Slide 3: hypothesis testing bullet ---
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2)
$tr3 = $body3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(3, 1)
$para3.Text = "placeholder"
$para3.Text = "Hypothesis testing: Once a control was established, I ran a series of hypothesis tests for types of crime associated with piracy with 𝛼 = 0.05:"

# --- 2. Slide 4: Portland legalized bullet ---
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2)
$tr4 = $body4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(2, 1)
$para4.Text = "placeholder"
$para4.Text = "We observe the same change in 2014 when Portland legalized, but again this is not statistically significant."

# --- 3. Slide 8: title split into two runs ---
$s8 = $p.Slides.Item(8)
$title8 = $s8.Shapes.Item(1)
$trTitle8 = $title8.TextFrame.TextRange
$firstPart = $trTitle8.Characters(1, 8)
$firstPart.Text = "Results "

# --- 4. Slide 8: conclusion bullet ---
$body8 = $s8.Shapes.Item(2)
$tr8 = $body8.TextFrame.TextRange
$para8 = $tr8.Paragraphs(1, 1)
$para8.Text = "placeholder"
$para8.Text = "Based on this data we cannot conclude that Weed Pirates are causing a significant increase in crime.  Instead, I observe a continuation of trends existing before legalization."
